# ==========================================================================
# Applies the "updates w/ ram comments" commit:
#  1. Inserts a new "timecalc" worksheet between "Sheet3" and "Sheet4"
#     with stride-time calculations for each country.
#  2. On "Sheet3": normalises the J/K "snapped" goal coordinates to a
#     constant (550, 300), tweaks a couple of the raw I column inputs,
#     and adds new A/B helper columns (linex/liney) that hold the old
#     (pre-snap) J/K values, row-shifted down by one.
# ==========================================================================

$wb = $excel.ActiveWorkbook

$sheet3 = $wb.Worksheets.Item("Sheet3")
$sheet4 = $wb.Worksheets.Item("Sheet4")

# --------------------------------------------------------------------
# 1. New "timecalc" worksheet, inserted right before "Sheet4"
# --------------------------------------------------------------------
$tc = $wb.Worksheets.Add($sheet4)
$tc.Name = "timecalc"

# ---- header row (B1:D1 reuse pre-existing shared strings) -------------
$sheet3.Range("A1").Copy() | Out-Null
$tc.Range("B1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats (border style 10)
$tc.Range("B1").Value = "Country"

$sheet3.Range("C1").Copy() | Out-Null
$tc.Range("C1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats (border style 10)
$tc.Range("C1").Value = "Ht"

$tc.Range("D1").Value = "stride length"

# The brand-new shared strings (502-510) must be *created* in the same
# order the original author typed them, so they land on the same indices:
#   E1, E10, E11, F1, F14, G14, F15, G1, H1
$tc.Range("E1").Value = "distance travelled in 1 stride (m)"      # -> 502
$tc.Range("E10").Value = "Assumption "                            # -> 503
$tc.Range("E11").Value = "no of strides per sec"                  # -> 504
$tc.Range("F1").Value = "Distance covered in 1 sec(m)"            # -> 505
$tc.Range("F14").Value = "4.968 m"                                # -> 506
$tc.Range("G14").Value = "1 sec"                                  # -> 507
$tc.Range("F15").Value = "1 m"                                    # -> 508
$tc.Range("G1").Value = "1 m is covered in(s) "                   # -> 509
$tc.Range("H1").Value = "27.5 m is covered in(s)"                 # -> 510

$tc.Range("F11").Value = 2
$tc.Range("G15").Formula = "=1/4.968"

# ---- per-country rows (2-7) --------------------------------------------
$countries = @(
  @{ Row = 2; Name = "germany";     Ht = 184;           D = 248.4 },
  @{ Row = 3; Name = "netherlands"; Ht = 183.7777778;    D = 248.10000003000002 },
  @{ Row = 4; Name = "belgium";     Ht = 180.83333329999999; D = 244.12499995499999 },
  @{ Row = 5; Name = "australia";   Ht = 180.7222222;    D = 243.97499997000003 },
  @{ Row = 6; Name = "spain";       Ht = 178.66666670000001; D = 241.20000004500002 },
  @{ Row = 7; Name = "india";       Ht = 177.05555559999999; D = 239.02500006 }
)

foreach ($row in $countries) {
  $r = $row.Row

  $sheet3.Range("A" + $r).Copy() | Out-Null
  $tc.Range("B" + $r).PasteSpecial(-4122) | Out-Null   # xlPasteFormats (border style 11/13)
  $tc.Range("B" + $r).Value = $row.Name

  $sheet3.Range("B" + $r).Copy() | Out-Null
  $tc.Range("C" + $r).PasteSpecial(-4122) | Out-Null   # xlPasteFormats (border style 12)
  $tc.Range("C" + $r).Value = $row.Ht

  $tc.Range("D" + $r).Value = $row.D
  $tc.Range("E" + $r).Formula = "=D" + $r + "/100"
  $tc.Range("F" + $r).Formula = "=E" + $r + "*2"
  $tc.Range("G" + $r).Formula = "=1/F" + $r
  $tc.Range("H" + $r).Formula = "=G" + $r + "*27.5"
}

# --------------------------------------------------------------------
# 2. "Sheet3" edits
# --------------------------------------------------------------------

# New A11/B11 mini-header ("linex"/"liney"), copied down from J1/K1.
$sheet3.Range("J1").Copy() | Out-Null
$sheet3.Range("A11").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$sheet3.Range("A11").Value = "linex"

$sheet3.Range("K1").Copy() | Out-Null
$sheet3.Range("B11").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$sheet3.Range("B11").Value = "liney"

# A12:B17 keep the *old* (pre-snap) J/K values, shifted down one row.
$oldJK = @(
  @{ Row = 12; J = 547.89712343688052; K = 305.08211670077372 },
  @{ Row = 13; J = 546.03410925962544; K = 303.96941790747411 },
  @{ Row = 14; J = 543.45725987381422; K = 302.71406755720756 },
  @{ Row = 15; J = 542.86111236372551; K = 300.00423088788011 },
  @{ Row = 16; J = 542.4529017787662;  K = 296.87978319055577 },
  @{ Row = 17; J = 543.65098258785088; K = 293.66038254582435 }
)
foreach ($row in $oldJK) {
  $r = $row.Row
  $sheet3.Range("A" + $r).Value = $row.J
  $sheet3.Range("B" + $r).Value = $row.K
}

# J2:K7 now snap exactly to the fixed goal point (550, 300).
for ($r = 2; $r -le 7; $r++) {
  $sheet3.Range("J" + $r).Value = 550
  $sheet3.Range("K" + $r).Value = 300
}

# A couple of raw inputs also changed.
$sheet3.Range("I3").Value = 290
$sheet3.Range("I6").Value = 310

# Column N auto-fit (header text unchanged, but width now cached).
$sheet3.Columns.Item(14).AutoFit() | Out-Null

# Selection / active range, per the diff.
$sheet3.Range("J3:K7").Select() | Out-Null
$tc.Range("H2").Select() | Out-Null
$sheet3.Activate() | Out-Null
